$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.156833648681641
$ws.Range("B1").Value = 2.370877027511597
$ws.Range("D1").Value = 2.400433540344238
$ws.Range("E1").Value = 1.223396420478821
